$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.878.57"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "2.426.33"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'554.95"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").Value = "'161.51"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.510"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("D9").Value = "'0.159"
$ws.Range("E9").Value = "  +7.62%  "
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "'4.78"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").Value = "'0.326"
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("D13").Value = "67.711.96"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "'0.0000169"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("D15").Value = "'23.11"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "'10.34"
$ws.Range("E16").Value = "  -3.83%  "
$ws.Range("D17").Value = "'335.84"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "'6.85"
$ws.Range("E18").Value = "  -1.97%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "'1.89"
$ws.Range("E20").Value = "  +1.79%  "
$ws.Range("D22").Value = "'66.62"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").Value = "'3.63"
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("D24").Value = "'8.08"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").Value = "0.0₃0812"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").Value = "'7.12"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "'422.97"
$ws.Range("E28").Value = "  -1.79%  "
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("D30").Value = "'1.60"
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("D31").Value = "'160.85"
$ws.Range("E31").Value = "  +2.70%  "
$ws.Range("D32").Value = "'18.94"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").Value = "'17.76"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  -5.48%  "
$ws.Range("E36").Value = "  -2.70%  "
$ws.Range("D37").Value = "'4.27"
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("D38").Value = "'1.47"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").Value = "'1.06"
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("D40").Value = "'2.02"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'3.34"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "'129.11"
$ws.Range("E42").Value = "  -1.87%  "
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").Value = "'0.478"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "'0.556"
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("D46").Value = "'0.0914"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("E48").Value = "  -5.37%  "
$ws.Range("D49").Value = "'16.64"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").Value = "0.0₆0205"
$ws.Range("E50").Value = "  +4.21%  "
$ws.Range("D51").Value = "'4.78"
$ws.Range("E51").Value = "  -6.20%  "
